$d = $word.ActiveDocument

function Insert-RawXml($range, [string]$bodyXml) {
    $pkg = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?>' +
           '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
           $bodyXml +
           '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $range.InsertXML($pkg)
}

# ---------------------------------------------------------------------------
# Edit 1: insert a new "Meta description" paragraph right after the title
# (Heading1) paragraph - bold "Meta description" run followed by a plain
# run with the rest of the sentence. A throw-away marker paragraph is
# appended so the insertion creates a real paragraph break instead of
# merging into the following "Gameplay" paragraph; the marker text is then
# deleted, which also removes the scratch run/paragraph mark it introduced.
# ---------------------------------------------------------------------------
$title = $d.Paragraphs.Item(1)
$afterTitle = $d.Range($title.Range.End, $title.Range.End)

$metaBody = '<w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r>' +
            '<w:r><w:t>: Read our review of Dragon Spark, an online slot game with beautiful visuals, multiple bonus features, and up to 16,807 ways to win. Play for free.</w:t></w:r></w:p>' +
            '<w:p><w:r><w:t>ZZZ_TMP_MARKER_ZZZ</w:t></w:r></w:p>'
Insert-RawXml $afterTitle $metaBody

$markerRange = $d.Content
$markerRange.Find.Execute("ZZZ_TMP_MARKER_ZZZ", $false, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

# Remember where the freshly-inserted meta-description paragraph ends, so
# later searches for its (duplicated) sentence text can be scoped to only
# look *after* this point and hit the real target near the end of the doc.
$metaParaEnd = $d.Paragraphs.Item(2).Range.End

# ---------------------------------------------------------------------------
# Edit 2: remove the duplicated bold title paragraph near the end of the
# document (the one that used to sit right before the italic meta-
# description paragraph).
# ---------------------------------------------------------------------------
$dupIndex = -1
for ($i = 2; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Play Dragon Spark for Free - Review of the Dragon-Themed Online Slot*") {
        $dupIndex = $i
    }
}
if ($dupIndex -ge 1) {
    $d.Paragraphs.Item($dupIndex).Range.Delete() | Out-Null
}

# ---------------------------------------------------------------------------
# Edit 3: replace the text of the trailing italic paragraph (old meta-
# description copy) with the new AI image-prompt text, keeping the italic
# run formatting intact. The search range starts right after the new
# meta-description paragraph, so the (now duplicated) sentence that lives
# inside it does not shadow the real target near the end of the document.
# ---------------------------------------------------------------------------
$docEnd = $d.Content.End
$findRange = $d.Range($metaParaEnd, $docEnd)
$findRange.Find.Execute("Read our review of Dragon Spark, an online slot game with beautiful visuals, multiple bonus features, and up to 16,807 ways to win. Play for free.", `
    $true, $false, $false, $false, $false, $true, 0, $false, "", 0) | Out-Null
$oldItalic = $d.Range($findRange.Start, $findRange.End)

$newItalicBody = '<w:p><w:r><w:rPr><w:i/></w:rPr><w:t>' +
    "Create a cartoon-style feature image for Dragon Spark that features a happy Maya warrior with glasses. The image should be eye-catching and showcase the theme of the game, including the dragon and princess elements. Please make sure to include the game logo somewhere in the image and use bright colors to grab the viewer's attention." +
    '</w:t></w:r></w:p>'
Insert-RawXml $oldItalic $newItalicBody

Write-Output "Edits applied. Paragraph count: $($d.Paragraphs.Count)"
